$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update r2 values for existing rows (A2:A8 / 2,3,4,5,10,15,20)
$ws.Range("B2").Value = 0.8279476871118658
$ws.Range("B3").Value = 0.8281712605073157
$ws.Range("B4").Value = 0.8278746192363403
$ws.Range("B5").Value = 0.8282773431781433
$ws.Range("B6").Value = 0.8300777276889197
$ws.Range("B7").Value = 0.8332427783193037
$ws.Range("B8").Value = 0.8306366459920242

# Remove the row for Embedding Size = 30 (row 9)
$ws.Range("A9:B9").Delete() | Out-Null
